# Import input xlsx format correctly for 1 sheet
# - fill in the previously-empty B4 cell
# - change the row-1 column-B header to "Đường đi"
# - append a second "Section" block (rows 5-8), mirroring rows 1-4,
#   with the header cell in column B reusing the original
#   "Section Name" text
# - leave the active selection on B1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete the existing block: B4 was empty, give it the same "1" used
# by B2/B3.
$ws.Range("B4").Value = 1
$ws.Rows.Item(4).RowHeight = 12.1

# The column-B header in row 1 becomes "Đường đi" ...
$ws.Range("B1").Value = "Đường đi"

# ... and a new block (rows 5-8) repeats the Section/1.1/1.2/1.3 layout
# of rows 1-4, reusing the original "Section Name" header text.
$ws.Range("A5").Value = "Section"
$ws.Range("B5").Value = "Section Name"
$ws.Range("A6").Value = 1.1
$ws.Range("B6").Value = 1
$ws.Range("A7").Value = 1.2
$ws.Range("B7").Value = 1
$ws.Range("A8").Value = 1.3
$ws.Range("B8").Value = 1

# Match the row heights already used on rows 1-4 for the new rows.
$ws.Rows.Item(5).RowHeight = 12.1
$ws.Rows.Item(6).RowHeight = 12.1
$ws.Rows.Item(7).RowHeight = 12.1
$ws.Rows.Item(8).RowHeight = 12.1

# Also mark rows 1-3 as explicitly sized (row 4 already covered above).
$ws.Rows.Item(1).RowHeight = 12.1
$ws.Rows.Item(2).RowHeight = 12.1
$ws.Rows.Item(3).RowHeight = 12.1

# Move the active cell/selection back up to B1.
$ws.Range("B1").Select()
